$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1008.9167
$ws.Range("I28").Value = 879.7059
$ws.Range("J28").Value = 1322.7142
$ws.Range("K28").Value = 879.7059
$ws.Range("L28").Value = 1322.7142
$ws.Range("M28").Value = -394.7059
$ws.Range("N28").Value = -2292.7142
$ws.Range("H32").Value = 3823.75
$ws.Range("I32").Value = 4131.6665
$ws.Range("J32").Value = 2900
$ws.Range("K32").Value = 4131.6665
$ws.Range("L32").Value = 2900
$ws.Range("M32").Value = -3805.6665
$ws.Range("N32").Value = -3552
$ws.Range("H41").Value = 274.66666
$ws.Range("I41").Value = 305.25
$ws.Range("J41").Value = 250.2
$ws.Range("K41").Value = 305.25
$ws.Range("L41").Value = 250.2
$ws.Range("M41").Value = 134.75
$ws.Range("N41").Value = -1130.2
$ws.Range("H51").Value = 2827.7144
$ws.Range("I51").Value = 2800
$ws.Range("J51").Value = 2829.8462
$ws.Range("K51").Value = 2800
$ws.Range("L51").Value = 2829.8462
$ws.Range("M51").Value = -2316
$ws.Range("N51").Value = -3797.8462
$ws.Range("H98").Value = 3243.3
$ws.Range("I98").Value = 3308.7368
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 3308.7368
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -1810.7368
$ws.Range("N98").Value = -4996
$ws.Range("H112").Value = 58824756
$ws.Range("I112").Value = 607.5
$ws.Range("J112").Value = 76924500
$ws.Range("K112").Value = 1822.5
$ws.Range("L112").Value = 230773500
$ws.Range("M112").Value = -714.5
$ws.Range("N112").Value = -230775716
$ws.Range("H122").Value = 3243.3
$ws.Range("I122").Value = 3308.7368
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 9926.2104
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -7476.2104
$ws.Range("N122").Value = -10900
$ws.Range("H129").Value = 1542.9615
$ws.Range("J129").Value = 1564.68
$ws.Range("L129").Value = 4694.04
$ws.Range("N129").Value = -14694.04

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3698.1667
$ws.Range("I61").Value = 2068.75
$ws.Range("J61").Value = 6957
$ws.Range("K61").Value = 2068.75
$ws.Range("L61").Value = 6957
$ws.Range("M61").Value = -1856.75
$ws.Range("N61").Value = -7381
$ws.Range("H74").Value = 3224.62
$ws.Range("I74").Value = 897.2105
$ws.Range("J74").Value = 4651.0967
$ws.Range("K74").Value = 897.2105
$ws.Range("L74").Value = 4651.0967
$ws.Range("M74").Value = -23.21050000000002
$ws.Range("N74").Value = -6399.0967
$ws.Range("H77").Value = 3224.62
$ws.Range("I77").Value = 897.2105
$ws.Range("J77").Value = 4651.0967
$ws.Range("K77").Value = 4486.0525
$ws.Range("L77").Value = 23255.4835
$ws.Range("M77").Value = -118.0524999999998
$ws.Range("N77").Value = -31991.4835
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H106").Value = 36000
$ws.Range("J106").Value = 36000
$ws.Range("L106").Value = 36000
$ws.Range("N106").Value = -38524
$ws.Range("H136").Value = 3698.1667
$ws.Range("I136").Value = 2068.75
$ws.Range("J136").Value = 6957
$ws.Range("K136").Value = 6206.25
$ws.Range("L136").Value = 20871
$ws.Range("M136").Value = -3656.25
$ws.Range("N136").Value = -25971

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 3022
$ws.Range("I54").Value = 3022
$ws.Range("K54").Value = 3022
$ws.Range("M54").Value = -2538
$ws.Range("H88").Value = 51628
$ws.Range("J88").Value = 51628
$ws.Range("L88").Value = 51628
$ws.Range("N88").Value = -52440
$ws.Range("H91").Value = 51628
$ws.Range("J91").Value = 51628
$ws.Range("L91").Value = 51628
$ws.Range("N91").Value = -54436

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1004.2727
$ws.Range("I35").Value = 1004.2727
$ws.Range("K35").Value = 1004.2727
$ws.Range("M35").Value = -710.2727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1226.1364
$ws.Range("I5").Value = 1443.8462
$ws.Range("J5").Value = 911.6667
$ws.Range("K5").Value = 4331.5386
$ws.Range("L5").Value = 2735.0001
$ws.Range("M5").Value = -4219.5386
$ws.Range("N5").Value = -2959.0001
$ws.Range("H22").Value = 1031.25
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 1046.6666
$ws.Range("K22").Value = 2400
$ws.Range("L22").Value = 3139.9998
$ws.Range("M22").Value = -2231
$ws.Range("N22").Value = -3477.9998
$ws.Range("H27").Value = 1031.25
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 1046.6666
$ws.Range("K27").Value = 2400
$ws.Range("L27").Value = 3139.9998
$ws.Range("M27").Value = -2298
$ws.Range("N27").Value = -3343.9998
$ws.Range("H43").Value = 3000
$ws.Range("J43").Value = 3000
$ws.Range("L43").Value = 9000
$ws.Range("N43").Value = -9228
$ws.Range("H56").Value = 2096.6667
$ws.Range("I56").Value = 2096.6667
$ws.Range("K56").Value = 2096.6667
$ws.Range("M56").Value = -1566.6667
$ws.Range("H94").Value = 3050.2666
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 3050.2666
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 9150.799800000001
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -10502.7998
$ws.Range("H96").Value = 47154736
$ws.Range("J96").Value = 47154736
$ws.Range("L96").Value = 141464208
$ws.Range("N96").Value = -141468326
$ws.Range("H100").Value = 2745
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 9000
$ws.Range("N100").Value = -10622
$ws.Range("H103").Value = 4858300
$ws.Range("I103").Value = 8500400
$ws.Range("J103").Value = 2166.6667
$ws.Range("K103").Value = 25501200
$ws.Range("L103").Value = 6500.000100000001
$ws.Range("M103").Value = -25500321
$ws.Range("N103").Value = -8258.000100000001
$ws.Range("H105").Value = 181602460
$ws.Range("J105").Value = 181602460
$ws.Range("L105").Value = 544807380
$ws.Range("N105").Value = -544812622
$ws.Range("H110").Value = 2890.9092
$ws.Range("J110").Value = 3683.3333
$ws.Range("L110").Value = 11049.9999
$ws.Range("N110").Value = -19229.9999
$ws.Range("H114").Value = 1954.5217
$ws.Range("I114").Value = 918.7
$ws.Range("J114").Value = 2751.3076
$ws.Range("K114").Value = 2756.1
$ws.Range("L114").Value = 8253.9228
$ws.Range("M114").Value = 497.8999999999996
$ws.Range("N114").Value = -14761.9228
$ws.Range("H122").Value = 11112255
$ws.Range("I122").Value = 17544338
$ws.Range("J122").Value = 2292.0908
$ws.Range("K122").Value = 157899042
$ws.Range("L122").Value = 20628.8172
$ws.Range("M122").Value = -157896592
$ws.Range("N122").Value = -25528.8172
$ws.Range("H129").Value = 1094.4615
$ws.Range("I129").Value = 594.75
$ws.Range("J129").Value = 1316.5555
$ws.Range("K129").Value = 1784.25
$ws.Range("L129").Value = 3949.6665
$ws.Range("M129").Value = 3215.75
$ws.Range("N129").Value = -13949.6665
$ws.Range("H131").Value = 1697595.9
$ws.Range("J131").Value = 1888610.9
$ws.Range("L131").Value = 5665832.699999999
$ws.Range("N131").Value = -5675912.699999999
$ws.Range("H135").Value = 1226.1364
$ws.Range("I135").Value = 1443.8462
$ws.Range("J135").Value = 911.6667
$ws.Range("K135").Value = 12994.6158
$ws.Range("L135").Value = 8205.0003
$ws.Range("M135").Value = -10459.6158
$ws.Range("N135").Value = -13275.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 2772.2
$ws.Range("J19").Value = 3095.238
$ws.Range("L19").Value = 3095.238
$ws.Range("N19").Value = -3671.238
$ws.Range("H122").Value = 7461.8
$ws.Range("I122").Value = 9391.9
$ws.Range("K122").Value = 28175.7
$ws.Range("M122").Value = -25725.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 460.5357
$ws.Range("I107").Value = 290.88235
$ws.Range("J107").Value = 722.7273
$ws.Range("K107").Value = 872.6470499999999
$ws.Range("L107").Value = 2168.1819
$ws.Range("M107").Value = 1047.35295
$ws.Range("N107").Value = -6008.1819
$ws.Range("H113").Value = 817.7778
$ws.Range("I113").Value = 656.4
$ws.Range("J113").Value = 1019.5
$ws.Range("K113").Value = 1969.2
$ws.Range("L113").Value = 3058.5
$ws.Range("M113").Value = 200.8000000000002
$ws.Range("N113").Value = -7398.5

Write-Host "Applied all changes"